# Daily attendance processing - 2025-12-18 07:31:53
#
# Applies the refreshed attendance-system data pull to the session
# analysis workbook:
#   - Swaps the "Recorded By" ordering from "System, <email>" to
#     "<email>, System" wherever both recorded the session.
#   - Refreshes the attendance-derived counts / percentages that moved
#     after the new data sync (Class Statistics + per-group stats).
#   - Marks the 5 "Not Recorded" SURGERY SEMINAR/SLIDE B1 sessions
#     (rows 39, 61, 210, 232, 254) as now Recorded, with their
#     attendance figures and row styling matching the rest of the
#     recorded rows.
#   - Narrows column I now that "Not Recorded" no longer needs the
#     extra width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write literal text into a cell without Excel's automatic
# number/percentage/date reinterpretation, and without forking a new
# cell style (keeps the destination's existing style index intact).
# Works by staging the literal text in a scratch cell via a text
# formula, then copying only the *value* onto the destination (which
# preserves the destination cell's own formatting).
# ---------------------------------------------------------------------
$scratch = $ws.Range("ZZ1000")
function Set-LiteralText($rangeAddress, $text) {
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)
}

# 1) Recorded By ordering swap: "System, <email>" -> "<email>, System"
$swapCells = @(
    "G3","G4","G24","G25","G46","G47","G68","G69","G90","G91","G111","G112",
    "G132","G133","G153","G154","G174","G175","G195","G196","G217","G218","G239","G240"
)
foreach ($addr in $swapCells) {
    $ws.Range($addr).Value = "dnasr281@gmail.com, System"
}

# 2) Class Statistics block (K4:L10)
$ws.Range("L6").Value = 78
$ws.Range("L7").Value = 0
Set-LiteralText "L9" "30.2%"
Set-LiteralText "L10" "81.6%"

# 3) Per-group statistics rows (K16:S26)
$ws.Range("O16").Value = 7
$ws.Range("P16").Value = 0
Set-LiteralText "R16" "31.8%"
Set-LiteralText "S16" "78.8%"

$ws.Range("O17").Value = 7
$ws.Range("P17").Value = 0
Set-LiteralText "R17" "31.8%"
Set-LiteralText "S17" "66.9%"

Set-LiteralText "S18" "87.8%"

$ws.Range("O24").Value = 7
$ws.Range("P24").Value = 0
Set-LiteralText "R24" "31.8%"
Set-LiteralText "S24" "71.4%"

$ws.Range("O25").Value = 7
$ws.Range("P25").Value = 0
Set-LiteralText "R25" "31.8%"
Set-LiteralText "S25" "77.8%"

$ws.Range("O26").Value = 7
$ws.Range("P26").Value = 0
Set-LiteralText "R26" "31.8%"
Set-LiteralText "S26" "72.9%"

# 4) Newly recorded sessions: re-style rows to the standard "Recorded"
#    look (same formatting as the other recorded data rows) and fill
#    in the attendance results that came back from the sync.
$newlyRecorded = @(
    @{ Row = 39;  Students = "23/31" },
    @{ Row = 61;  Students = "9/19" },
    @{ Row = 210; Students = "16/27" },
    @{ Row = 232; Students = "18/29" },
    @{ Row = 254; Students = "18/29" }
)

$ws.Range("A2:I2").Copy()
foreach ($item in $newlyRecorded) {
    $ws.Range("A" + $item.Row + ":I" + $item.Row).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

foreach ($item in $newlyRecorded) {
    $r = $item.Row
    $ws.Range("G" + $r).Value = "System"
    $ws.Range("H" + $r).Value = $item.Students
    $ws.Range("I" + $r).Value = "Recorded"
}

# 5) Single attendance-count correction outside the above rows
$ws.Range("H83").Value = "17/21"

# 6) Column I no longer needs to fit "Not Recorded"; narrow it back
#    down to match the other short status-style columns.
$ws.Columns.Item(9).ColumnWidth = (10 - 6/7)

# Clean up the scratch cell used for literal-text staging.
$scratch.Clear()
